$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header columns to support plans without sticky-ids or levels
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the active selection to A2
$ws.Range("A2").Select()
